$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.380.76'
$ws.Range("D3").Value = '1.678.79'
$ws.Range("E3").Value = '  +1.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5297'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("E8").Value = '  +3.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06445'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.27'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.09%  '
$ws.Range("D12").Value = '1.674.89'
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("E13").Value = '  +2.85%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5608'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.92%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '0.0₅8438'
$ws.Range("E15").Value = '  +6.35%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '26.393.28'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.835'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.48%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.16%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.81%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.396'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.28%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.27%  '
$ws.Range("B25").Value = 'Stellar'
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1264'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.10%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.486'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.75%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.50%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.435'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.41%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06200'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.03%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.276'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.79%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.552'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.22%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.464'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.99%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.703'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.56%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.019'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.10%  '
$ws.Range("B35").Value = 'MXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.789'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.94%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.404'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.02%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5747'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.37%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01644'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.22%  '
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.943'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8687'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("B41").Value = 'Maker'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D41").Value = '1.055.73'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9995'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.827.73'
$ws.Range("E44").Value = '  +1.41%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.08'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.03%  '
$ws.Range("D46").Value = '0.0₈108'
$ws.Range("E46").Value = '  +2.70%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.166'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.11%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.0000'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05199'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.053'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.19%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09959'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.41%  '
